# Update "想去人数" (number of people interested) figures that changed
# between data pulls, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5200
$ws1.Range("F6").Value = 302
$ws1.Range("F7").Value = 788
$ws1.Range("F8").Value = 281

# Sheet "全部类型" (all types combined list)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5200
$ws4.Range("F6").Value = 302
$ws4.Range("F7").Value = 788
$ws4.Range("F9").Value = 281
